# Refresh the cryptocurrency price/volume snapshot to the values captured
# by the latest GitHub Actions scrape run.
#   Columns: B = Coin name, C = Link, D = Price, E = Volume(1h)
#   Rows 9/10 (Cardano <-> OKB) and 29/30 (EthereumClassic <-> Kaspa) swapped
#   places in the ranking, so both name/link/price/volume cells move.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell on this sheet is stored as text (even the Price column, which
# looks numeric, e.g. '233.88' or '36.454.25'). Writing a numeric-looking
# string straight into .Value lets Excel silently re-type the cell as a
# number, so Set-TextValue forces text with a leading quote-prefix, just
# like typing '233.88 into a cell would.
function Set-TextValue($range, $text) {
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($range).Value = "'" + $text
    } else {
        $ws.Range($range).Value = $text
    }
}

Set-TextValue "D2" '36.454.25'
$ws.Range("E2").Value = '  -2.80%  '
Set-TextValue "D3" '1.984.47'
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("E4").Value = '  -0.14%  '
Set-TextValue "D5" '233.88'
$ws.Range("E5").Value = '  -10.88%  '
Set-TextValue "D6" '0.599'
$ws.Range("E6").Value = '  -3.57%  '
$ws.Range("E7").Value = '  -0.02%  '
Set-TextValue "D8" '54.32'
$ws.Range("E8").Value = '  -3.08%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue "D9" '0.372'
$ws.Range("E9").Value = '  -3.74%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D10" '58.58'
$ws.Range("E10").Value = '  +3.11%  '
Set-TextValue "D11" '0.0745'
$ws.Range("E11").Value = '  -3.81%  '
Set-TextValue "D12" '0.0984'
$ws.Range("E12").Value = '  -3.36%  '
Set-TextValue "D13" '14.19'
$ws.Range("E13").Value = '  -1.34%  '
Set-TextValue "D14" '2.273.50'
$ws.Range("E14").Value = '  -1.69%  '
Set-TextValue "D15" '19.83'
$ws.Range("E15").Value = '  -5.38%  '
Set-TextValue "D16" '0.755'
$ws.Range("E16").Value = '  -6.56%  '
Set-TextValue "D17" '5.04'
$ws.Range("E17").Value = '  -4.23%  '
Set-TextValue "D18" '1.989.98'
$ws.Range("E18").Value = '  -1.81%  '
Set-TextValue "D19" '36.357.76'
$ws.Range("E19").Value = '  -2.78%  '
Set-TextValue "D20" '67.57'
$ws.Range("E20").Value = '  -3.11%  '
Set-TextValue "D21" '0.0₃0801'
$ws.Range("E21").Value = '  -5.01%  '
Set-TextValue "D22" '5.28'
$ws.Range("E22").Value = '  +1.75%  '
Set-TextValue "D23" '221.30'
$ws.Range("E23").Value = '  -3.16%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  +0.86%  '
Set-TextValue "D26" '2.40'
$ws.Range("E26").Value = '  -11.19%  '
Set-TextValue "D27" '161.55'
$ws.Range("E27").Value = '  -2.10%  '
Set-TextValue "D28" '8.54'
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D29" '18.82'
$ws.Range("E29").Value = '  -4.64%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D30" '0.125'
$ws.Range("E30").Value = '  -2.28%  '
Set-TextValue "D31" '1.32'
$ws.Range("E31").Value = '  -1.35%  '
$ws.Range("E32").Value = '  -3.24%  '
Set-TextValue "D33" '4.37'
$ws.Range("E33").Value = '  -6.19%  '
Set-TextValue "D34" '0.0605'
$ws.Range("E34").Value = '  -6.93%  '
Set-TextValue "D35" '4.23'
$ws.Range("E35").Value = '  -6.65%  '
$ws.Range("E36").Value = '  -3.38%  '
Set-TextValue "D37" '0.999'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("E38").Value = '  -2.44%  '
Set-TextValue "D39" '3.21'
$ws.Range("E39").Value = '  -4.71%  '
Set-TextValue "D40" '5.40'
$ws.Range("E40").Value = '  +4.28%  '
$ws.Range("E41").Value = '  -1.09%  '
Set-TextValue "D42" '1.451.61'
$ws.Range("E42").Value = '  +3.75%  '
Set-TextValue "D43" '0.0915'
$ws.Range("E43").Value = '  -2.91%  '
Set-TextValue "D44" '0.0201'
$ws.Range("E44").Value = '  -6.21%  '
Set-TextValue "D45" '1.09'
$ws.Range("E45").Value = '  -11.48%  '
Set-TextValue "D46" '88.35'
$ws.Range("E46").Value = '  -2.26%  '
Set-TextValue "D47" '0.992'
$ws.Range("E47").Value = '  -3.50%  '
Set-TextValue "D48" '14.75'
$ws.Range("E48").Value = '  -6.10%  '
$ws.Range("E49").Value = '  -1.10%  '
Set-TextValue "D50" '6.75'
$ws.Range("E50").Value = '  -4.42%  '
Set-TextValue "D51" '2.166.82'
$ws.Range("E51").Value = '  -1.68%  '
